$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 205, pushing the existing rows 205..242 down to 206..243.
$ws.Rows.Item(205).Insert()

# Populate the newly inserted row 205 with a new weekly record: same as the
# record that is now in row 206 (the former row 205), except for a new date.
$ws.Range("A206:T206").Copy()
$ws.Range("A205:T205").PasteSpecial()

$ws.Cells.Item(205, 4).Value = "5/13/2022"
